# Apply the commit: "Updated results (single bin), and moving the PCA
# reduction stuff in."
#
# 1. Fill in the new (single-bin) results row on the SEMAINE sheet, row 4,
#    columns B:M - highlighting every column except the Pr/Rec pair (E, I)
#    with the built-in "Good" (green) cell style, matching the existing
#    highlighting convention used elsewhere in the sheet.
# 2. Update the selection caret on each sheet and switch the active tab
#    from BP4D to SEMAINE (the sheet the author was actively working in).

$wb = $excel.ActiveWorkbook

$bp4d    = $wb.Worksheets.Item("BP4D")
$semaine = $wb.Worksheets.Item("SEMAINE")

# --- New single-bin results row (SEMAINE!B4:M4) -----------------------------
$values = [ordered]@{
    "B4" = 0.46899999999999997
    "C4" = 0.51100000000000001
    "D4" = 0.48899999999999999
    "E4" = 0.51600000000000001
    "F4" = 0.46300000000000002
    "G4" = 0.48799999999999999
    "H4" = 0.23899999999999999
    "I4" = 0.61199999999999999
    "J4" = 0.34300000000000003
    "K4" = 0.28999999999999998
    "L4" = 0.54600000000000004
    "M4" = 0.379
}

# Columns E (Pr) and I (Pr) are left with the default style; every other
# column gets the "Good" highlight, same as the rest of the workbook.
$highlighted = @("B4", "C4", "D4", "F4", "G4", "H4", "J4", "K4", "L4", "M4")

foreach ($addr in $values.Keys) {
    $cell = $semaine.Range($addr)
    $cell.Value = $values[$addr]
}

foreach ($addr in $highlighted) {
    $semaine.Range($addr).Style = "Good"
}

# --- Selection / active tab --------------------------------------------------
# Set BP4D's lingering selection first ...
$bp4d.Range("D4").Select() | Out-Null
# ... then finish on SEMAINE so it ends up as the active/selected tab.
$semaine.Range("L8").Select() | Out-Null
